$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 135.42857
$ws.Range("I12").Value = 139.6
$ws.Range("K12").Value = 139.6
$ws.Range("M12").Value = 30.40000000000001
# Row 32
$ws.Range("H32").Value = 2348
$ws.Range("I32").Value = 294
$ws.Range("K32").Value = 294
$ws.Range("M32").Value = 32
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
# Row 51
$ws.Range("H51").Value = 8538.462
$ws.Range("J51").Value = 9041.666999999999
$ws.Range("L51").Value = 9041.666999999999
$ws.Range("N51").Value = -10009.667
# Row 53
$ws.Range("H53").Value = 256
$ws.Range("I53").Value = 272.6
$ws.Range("K53").Value = 272.6
$ws.Range("M53").Value = 364.4
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
# Row 62
$ws.Range("H62").Value = 10999.333
# Row 65
$ws.Range("H65").Value = 10999.333
# Row 70
$ws.Range("H70").Value = 9161.666999999999
$ws.Range("I70").Value = 6333.3335
$ws.Range("K70").Value = 19000.0005
$ws.Range("M70").Value = -18730.0005
# Row 73
$ws.Range("H73").Value = 9161.666999999999
$ws.Range("I73").Value = 6333.3335
$ws.Range("K73").Value = 19000.0005
$ws.Range("M73").Value = -18064.0005
# Row 80
$ws.Range("H80").Value = 1401
$ws.Range("I80").Value = 2334.6667
$ws.Range("K80").Value = 7004.000100000001
$ws.Range("M80").Value = -6006.000100000001
# Row 83
$ws.Range("H83").Value = 1401
$ws.Range("I83").Value = 2334.6667
$ws.Range("K83").Value = 21012.0003
$ws.Range("M83").Value = -16020.0003
# Row 88
$ws.Range("H88").Value = 4161.5557
$ws.Range("J88").Value = 5590.4
$ws.Range("L88").Value = 5590.4
$ws.Range("N88").Value = -6402.4
# Row 91
$ws.Range("H91").Value = 4161.5557
$ws.Range("J91").Value = 5590.4
$ws.Range("L91").Value = 5590.4
$ws.Range("N91").Value = -8398.4
# Row 132
$ws.Range("H132").Value = 4281.75
$ws.Range("I132").Value = 4281.75
$ws.Range("K132").Value = 12845.25
$ws.Range("M132").Value = -10315.25
# Row 141
$ws.Range("H141").Value = 1097
$ws.Range("I141").Value = 1097
$ws.Range("K141").Value = 3291
$ws.Range("M141").Value = 1889

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 9140
$ws.Range("I61").Value = 10364
$ws.Range("J61").Value = 8375
$ws.Range("K61").Value = 10364
$ws.Range("L61").Value = 8375
$ws.Range("M61").Value = -10152
$ws.Range("N61").Value = -8799
# Row 135
$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140
# Row 136
$ws.Range("H136").Value = 9140
$ws.Range("I136").Value = 10364
$ws.Range("J136").Value = 8375
$ws.Range("K136").Value = 31092
$ws.Range("L136").Value = 25125
$ws.Range("M136").Value = -28542
$ws.Range("N136").Value = -30225

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5554.6
$ws.Range("I20").Value = 4468.25
$ws.Range("J20").Value = 9900
$ws.Range("K20").Value = 4468.25
$ws.Range("L20").Value = 9900
$ws.Range("M20").Value = -4221.25
$ws.Range("N20").Value = -10394
# Row 94
$ws.Range("H94").Value = 2604.5
$ws.Range("I94").Value = 2604.5
$ws.Range("K94").Value = 2604.5
$ws.Range("M94").Value = -2153.5
# Row 105
$ws.Range("H105").Value = 12998
$ws.Range("I105").Value = 12998
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 12998
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -11251
$ws.Range("N105").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 71429610
$ws.Range("I4").Value = 1116.1538
$ws.Range("K4").Value = 3348.4614
$ws.Range("M4").Value = -3236.4614
# Row 11
$ws.Range("H11").Value = 1446.4286
$ws.Range("I11").Value = 576.6667
$ws.Range("K11").Value = 1730.0001
$ws.Range("M11").Value = -1590.0001
# Row 81
$ws.Range("H81").Value = 2999
$ws.Range("J81").Value = 2999
$ws.Range("L81").Value = 8997
$ws.Range("N81").Value = -11243
# Row 84
$ws.Range("H84").Value = 2999
$ws.Range("J84").Value = 2999
$ws.Range("L84").Value = 26991
$ws.Range("N84").Value = -38223
# Row 97
$ws.Range("H97").Value = 400
$ws.Range("I97").Value = 150
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 450
$ws.Range("L97").Value = 2700
$ws.Range("M97").Value = 46
$ws.Range("N97").Value = -3692

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 23.05
$ws.Range("I2").Value = 25.733334
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 25.733334
$ws.Range("L2").Value = 15
$ws.Range("M2").Value = 87.266666
$ws.Range("N2").Value = -241
# Row 80
$ws.Range("H80").Value = 5665.6665
$ws.Range("I80").Value = 7499
$ws.Range("J80").Value = 1999
$ws.Range("K80").Value = 7499
$ws.Range("L80").Value = 1999
$ws.Range("M80").Value = -6501
$ws.Range("N80").Value = -3995
# Row 83
$ws.Range("H83").Value = 5665.6665
$ws.Range("I83").Value = 7499
$ws.Range("J83").Value = 1999
$ws.Range("K83").Value = 37495
$ws.Range("L83").Value = 9995
$ws.Range("M83").Value = -32503
$ws.Range("N83").Value = -19979
# Row 107
$ws.Range("H107").Value = 885.5714
$ws.Range("I107").Value = 866.5
$ws.Range("K107").Value = 866.5
$ws.Range("M107").Value = 1053.5
# Row 122
$ws.Range("H122").Value = 899.6667
$ws.Range("I122").Value = 899.6667
$ws.Range("K122").Value = 2699.0001
$ws.Range("M122").Value = -249.0001000000002
# Row 132
$ws.Range("H132").Value = 6165
$ws.Range("J132").Value = 8149.25
$ws.Range("L132").Value = 24447.75
$ws.Range("N132").Value = -29507.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3875
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 133
$ws.Range("H133").Value = 26000
$ws.Range("J133").Value = 26000
$ws.Range("L133").Value = 26000
$ws.Range("N133").Value = -31060
# Row 14
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 81
$ws.Range("H81").Value = 5566.5
$ws.Range("I81").Value = 1966.3334
$ws.Range("K81").Value = 3932.6668
$ws.Range("M81").Value = -2871.6668
# Row 84
$ws.Range("H84").Value = 5566.5
$ws.Range("I84").Value = 1966.3334
$ws.Range("K84").Value = 19663.334
$ws.Range("M84").Value = -14359.334
# Row 96
$ws.Range("H96").Value = 971.5
$ws.Range("I96").Value = 971.5
$ws.Range("K96").Value = 971.5
$ws.Range("M96").Value = 401.5
